# Update the Btc-Erbb2 LR-pairs sheet with the new TPM-derived values.
#
# The sending/target cluster grid grows from the original FAPs/MuSCs x
# {ECs,FAPs,MuSCs} set to a full FAPs/MuSCs x {ECs,FAPs,Inflammatory-Mac,
# MuSCs,Resolving-Mac} grid (new target clusters "Inflammatory-Mac" and
# "Resolving-Mac" are introduced), and every metric column is recomputed
# with the new TPM numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each inner array is: row number, then the 20 column values (A..T) in
# order: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression
# value, Ligand total expression value, Ligand derived specificity of
# average/total expression value, Receptor-expressing cells, Receptor
# detection rate, Receptor average/total expression value, Receptor
# derived specificity of average/total expression value, Edge average/total
# expression weight, Edge average/total expression derived specificity.
$rows = @(
  @(2, "FAPs", "Btc", "Erbb2", "ECs", 3, 1, 0.637502, 1.912506, 0.5022738159769057, 0.5022738159769057, 3, 1, 3.248861, 9.746582999999999, 0.3767295193213461, 0.376729519321346, 2.071155385222, 18.640398466998, 0.1892213732606779, 0.1892213732606779),
  @(3, "FAPs", "Btc", "Erbb2", "FAPs", 3, 1, 0.637502, 1.912506, 0.5022738159769057, 0.5022738159769057, 3, 1, 3.599001333333334, 10.797004, 0.417330886838049, 0.4173308868380489, 2.294370548002667, 20.649334932024, 0.2096143770571731, 0.209614377057173),
  @(4, "FAPs", "Btc", "Erbb2", "Inflammatory-Mac", 3, 1, 0.637502, 1.912506, 0.5022738159769057, 0.5022738159769057, 1, 0.3333333333333333, 0.02714433333333334, 0.08143300000000001, 0.003147586692371591, 0.003147586692371591, 0.01730456678866667, 0.155741101098, 0.001580950379095606, 0.001580950379095606),
  @(5, "FAPs", "Btc", "Erbb2", "MuSCs", 3, 1, 0.637502, 1.912506, 0.5022738159769057, 0.5022738159769057, 3, 1, 1.711959, 5.135877, 0.1985143381535413, 0.1985143381535413, 1.091377286418, 9.822395577762, 0.09970855415050904, 0.09970855415050901),
  @(6, "FAPs", "Btc", "Erbb2", "Resolving-Mac", 3, 1, 0.637502, 1.912506, 0.5022738159769057, 0.5022738159769057, 2, 0.6666666666666666, 0.03689, 0.11067, 0.004277668994692127, 0.004277668994692126, 0.02351744878, 0.21165703902, 0.002148561129450109, 0.002148561129450108),
  @(7, "MuSCs", "Btc", "Erbb2", "ECs", 3, 1, 0.63173, 1.89519, 0.4977261840230943, 0.4977261840230943, 3, 1, 3.248861, 9.746582999999999, 0.3767295193213461, 0.376729519321346, 2.05240295953, 18.47162663577, 0.1875081460606682, 0.1875081460606681),
  @(8, "MuSCs", "Btc", "Erbb2", "FAPs", 3, 1, 0.63173, 1.89519, 0.4977261840230943, 0.4977261840230943, 3, 1, 3.599001333333334, 10.797004, 0.417330886838049, 0.4173308868380489, 2.273597112306667, 20.46237401076, 0.2077165097808759, 0.2077165097808759),
  @(9, "MuSCs", "Btc", "Erbb2", "Inflammatory-Mac", 3, 1, 0.63173, 1.89519, 0.4977261840230943, 0.4977261840230943, 1, 0.3333333333333333, 0.02714433333333334, 0.08143300000000001, 0.003147586692371591, 0.003147586692371591, 0.01714788969666667, 0.15433100727, 0.001566636313275985, 0.001566636313275985),
  @(10, "MuSCs", "Btc", "Erbb2", "MuSCs", 3, 1, 0.63173, 1.89519, 0.4977261840230943, 0.4977261840230943, 3, 1, 1.711959, 5.135877, 0.1985143381535413, 0.1985143381535413, 1.08149585907, 9.73346273163, 0.09880578400303228, 0.09880578400303226),
  @(11, "MuSCs", "Btc", "Erbb2", "Resolving-Mac", 3, 1, 0.63173, 1.89519, 0.4977261840230943, 0.4977261840230943, 2, 0.6666666666666666, 0.03689, 0.11067, 0.004277668994692127, 0.004277668994692126, 0.0233045197, 0.2097406773, 0.002129107865242019, 0.002129107865242018)
)

foreach ($row in $rows) {
  $r = $row[0]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item($r, $cols[$i]).Value = $row[$i + 1]
  }
}
